# Desenvolupant mòdul de comparació.
# Reorders the B..H content of rows 2-33 on the active sheet (each row's
# data row moves to a different row, while column A stays fixed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$finalData = @{
    2  = @{ B = "ATATATATAT"; D = 0; E = 0; F = 0; G = "";          H = ""  }
    3  = @{ B = "AAATATATAT"; D = 1; E = 0; F = 1; G = "1";         H = ""  }
    4  = @{ B = "ATTAATATAT"; D = 0; E = 1; F = 1; G = "";          H = "3" }
    5  = @{ B = "TAATATATAT"; D = 0; E = 1; F = 1; G = "";          H = "1" }
    6  = @{ B = "ATAAATATAT"; D = 1; E = 0; F = 1; G = "3";         H = ""  }
    7  = @{ B = "ATATATTAAT"; D = 0; E = 1; F = 1; G = "";          H = "7" }
    8  = @{ B = "ATATATAAAT"; D = 1; E = 0; F = 1; G = "7";         H = ""  }
    9  = @{ B = "ATATTAATAT"; D = 0; E = 1; F = 1; G = "";          H = "5" }
    10 = @{ B = "ATATAAATAT"; D = 1; E = 0; F = 1; G = "5";         H = ""  }
    11 = @{ B = "AATAATATAT"; D = 1; E = 1; F = 2; G = "1";         H = "3" }
    12 = @{ B = "ATAAATAAAT"; D = 2; E = 0; F = 2; G = "3, 7";      H = ""  }
    13 = @{ B = "AAATATTAAT"; D = 1; E = 1; F = 2; G = "1";         H = "7" }
    14 = @{ B = "AAATAAATAT"; D = 2; E = 0; F = 2; G = "1, 5";      H = ""  }
    15 = @{ B = "ATATAAAAAT"; D = 2; E = 0; F = 2; G = "5, 7";      H = ""  }
    16 = @{ B = "ATAATAATAT"; D = 1; E = 1; F = 2; G = "3";         H = "5" }
    17 = @{ B = "ATAAAAATAT"; D = 2; E = 0; F = 2; G = "3, 5";      H = ""  }
    18 = @{ B = "TAATATAAAT"; D = 1; E = 1; F = 2; G = "7";         H = "1" }
    19 = @{ B = "AAAAATATAT"; D = 2; E = 0; F = 2; G = "1, 3";      H = ""  }
    20 = @{ B = "ATATAATAAT"; D = 1; E = 1; F = 2; G = "5";         H = "7" }
    21 = @{ B = "AAATTAATAT"; D = 1; E = 1; F = 2; G = "1";         H = "5" }
    22 = @{ B = "ATTAATAAAT"; D = 1; E = 1; F = 2; G = "7";         H = "3" }
    23 = @{ B = "AAATATAAAT"; D = 2; E = 0; F = 2; G = "1, 7";      H = ""  }
    24 = @{ B = "TAAAATATAT"; D = 1; E = 1; F = 2; G = "3";         H = "1" }
    25 = @{ B = "ATAAAAAAAT"; D = 3; E = 0; F = 3; G = "3, 5, 7";   H = ""  }
    26 = @{ B = "AAAAAAATAT"; D = 3; E = 0; F = 3; G = "1, 3, 5";   H = ""  }
    27 = @{ B = "AAAATAATAT"; D = 2; E = 1; F = 3; G = "1, 3";      H = "5" }
    28 = @{ B = "AAATAATAAT"; D = 2; E = 1; F = 3; G = "1, 5";      H = "7" }
    29 = @{ B = "AAATAAAAAT"; D = 3; E = 0; F = 3; G = "1, 5, 7";   H = ""  }
    30 = @{ B = "AATAATAAAT"; D = 2; E = 1; F = 3; G = "1, 7";      H = "3" }
    31 = @{ B = "TAAAATAAAT"; D = 2; E = 1; F = 3; G = "3, 7";      H = "1" }
    32 = @{ B = "AAAAATAAAT"; D = 3; E = 0; F = 3; G = "1, 3, 7";   H = ""  }
    33 = @{ B = "AAAAAAAAAT"; D = 4; E = 0; F = 4; G = "1, 3, 5, 7"; H = "" }
}

foreach ($r in $finalData.Keys) {
    $row = $finalData[$r]
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
}
